# The workbook records one market-day observation per row for "Zanahoria"
# (carrots) at "Terminal Hortofrutícola Agro Chillán". A new, more recent
# observation is inserted as row 262, which pushes all the existing rows
# (previously 262-346) down by one (now 263-347). No other data changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 262 - this shifts rows 262:346 down to 263:347,
# carrying their existing values/styles with them.
$ws.Rows("262:262").Insert()

# Populate the newly inserted row 262 with the new observation.
$ws.Range("A262").Value = 7
$ws.Range("B262").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C262").Value = "Ñuble"
$ws.Range("D262").Value = 44876
$ws.Range("D262").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E262").Value = 16
$ws.Range("F262").Value = 100114013
$ws.Range("G262").Value = "Zanahoria"
$ws.Range("H262").Value = "Sin especificar"
$ws.Range("I262").Value = "Primera"
$ws.Range("J262").Value = 160
$ws.Range("K262").Value = 13000
$ws.Range("L262").Value = 14000
$ws.Range("M262").Value = 13500
$ws.Range("N262").Value = "`$/saco 20 kilos"
$ws.Range("O262").Value = "Región de Ñuble"
$ws.Range("P262").Value = 675
$ws.Range("Q262").Value = 20
$ws.Range("R262").Value = "Hortaliza"
